$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.200.96'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.906.60'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.35'
$ws.Range('E5').Value = '  +2.71%  '
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.77'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.366'
$ws.Range('E9').Value = '  +4.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.90'
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0757'
$ws.Range('E11').Value = '  +3.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0985'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '2.185.54'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '13.03'
$ws.Range('E14').Value = '  +5.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.748'
$ws.Range('E15').Value = '  +5.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.98'
$ws.Range('E16').Value = '  +2.62%  '
$ws.Range('D17').Value = '1.909.71'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '35.188.29'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.10'
$ws.Range('E19').Value = '  +1.94%  '
$ws.Range('D20').Value = '0.0₃0838'
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.11'
$ws.Range('E21').Value = '  +4.36%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '243.42'
$ws.Range('E22').Value = '  +0.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.07'
$ws.Range('E23').Value = '  +4.15%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  +6.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.34'
$ws.Range('E26').Value = '  -1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.04'
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.57'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('D31').Value = '4.128.63'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').Value = '  +14.32%  '
$ws.Range('E33').Value = '  +2.93%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0584'
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.61'
$ws.Range('E35').Value = '  +19.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.20'
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.864'
$ws.Range('E38').Value = '  -11.78%  '
$ws.Range('E39').Value = '  -1.89%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.47'
$ws.Range('E40').Value = '  +6.90%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.49'
$ws.Range('E41').Value = '  +9.55%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0217'
$ws.Range('E42').Value = '  +3.28%  '
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0651'
$ws.Range('E44').Value = '  -1.51%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.48'
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.339.21'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('E48').Value = '  -1.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.64'
$ws.Range('E49').Value = '  +0.94%  '
$ws.Range('B50').Value = 'Gas'
$ws.Range('C50').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.94'
$ws.Range('E50').Value = '  -7.42%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0752'
$ws.Range('E51').Value = '  +6.55%  '
